# Swap the "search" (TC2) and "cancel" (TC3) step contents.
# Before the edit:
#   B20/D20 (TC2, step #2) = "search" description/result
#   B28/D28 (TC3, step #2) = "cancel" description/result
# After the edit:
#   B20/D20 (TC2, step #2) = "cancel" description/result
#   B28/D28 (TC3, step #2) = "search" description/result

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$searchDesc   = "Chefe Indica alguns parâmetros específicos para a busca; Informa o nome do beneficiário; Filtra a listagem de solicitações."
$searchResult = "SYSTEM Exibe uma nova listagem de solicitações, de acordo com os filtros informados pelo usuário."
$cancelDesc   = "Chefe Clica para realizar o cancelamento de uma diária."
$cancelResult = "SYSTEM Verifica que a solicitação está em situação SOLICITADA; Exibe mensagem de confirmação (MSG987 - Cancelar solicitação de diária) para o usuário (que deve confirmar); Cancela a diária, mudando sua situação para CANCELADA (ver diagrama de estados da diária)."

$ws.Range("B20").Value = $cancelDesc
$ws.Range("D20").Value = $cancelResult

$ws.Range("B28").Value = $searchDesc
$ws.Range("D28").Value = $searchResult
